$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.922.16"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.631.85"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("E4").Value = "  +0.40%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "214.28"
$c.NumberFormat = "General"
$ws.Range("E5").Value = "  +0.77%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.517"
$c.NumberFormat = "General"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +0.27%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "28.36"
$c.NumberFormat = "General"
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("E9").Value = "  +0.80%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0607"
$c.NumberFormat = "General"
$ws.Range("E10").Value = "  +0.51%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0906"
$c.NumberFormat = "General"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "1.867.19"
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("D13").Value = "1.632.00"
$ws.Range("E13").Value = "  +1.93%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.561"
$c.NumberFormat = "General"
$ws.Range("E14").Value = "  +1.47%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "9.13"
$c.NumberFormat = "General"
$ws.Range("E15").Value = "  +13.21%  "
$ws.Range("D16").Value = "29.942.86"
$ws.Range("E16").Value = "  +0.62%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.84"
$c.NumberFormat = "General"
$ws.Range("E17").Value = "  +1.41%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "63.94"
$c.NumberFormat = "General"
$ws.Range("E18").Value = "  -0.23%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "241.41"
$c.NumberFormat = "General"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").Value = "0.0₃0700"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  +0.17%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.12"
$c.NumberFormat = "General"
$ws.Range("E22").Value = "  +1.89%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.75"
$c.NumberFormat = "General"
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("E24").Value = "  +2.83%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "159.44"
$c.NumberFormat = "General"
$ws.Range("E25").Value = "  +2.71%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "15.49"
$c.NumberFormat = "General"
$ws.Range("E26").Value = "  +0.04%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.109"
$c.NumberFormat = "General"
$ws.Range("E27").Value = "  +0.27%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "6.60"
$c.NumberFormat = "General"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("E29").Value = "  +0.37%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0486"
$c.NumberFormat = "General"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("E31").Value = "  +4.13%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.36"
$c.NumberFormat = "General"
$ws.Range("E32").Value = "  +3.79%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.17"
$c.NumberFormat = "General"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").Value = "1.424.07"
$ws.Range("E34").Value = "  -0.19%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.64"
$c.NumberFormat = "General"
$ws.Range("E35").Value = "  +4.50%  "
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("E39").Value = "  -0.30%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "75.45"
$c.NumberFormat = "General"
$ws.Range("E40").Value = "  +12.11%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.553"
$c.NumberFormat = "General"
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("E42").Value = "  +2.31%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.826"
$c.NumberFormat = "General"
$ws.Range("E43").Value = "  +0.91%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0491"
$c.NumberFormat = "General"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.NumberFormat = "General"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.02"
$c.NumberFormat = "General"
$ws.Range("E46").Value = "  +2.07%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "52.79"
$c.NumberFormat = "General"
$ws.Range("E47").Value = "  -4.62%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.773.99"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "5.34"
$c.NumberFormat = "General"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").Value = "0.0₆0114"
$ws.Range("E50").Value = "  +11.96%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "90.52"
$c.NumberFormat = "General"
$ws.Range("E51").Value = "  +4.39%  "

Write-Host "Applied 96 cell updates"
